$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as Text so Excel doesn't
# auto-convert purely numeric-looking strings (e.g. "250.90") into
# numbers and strip the meaningful trailing zero.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.157.59"
$ws.Range("E2").Value = "  -1.45%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.178.24"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "250.90"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -3.02%  "

# Row 7 - Solana
Set-TextValue "D7" "66.48"
$ws.Range("E7").Value = "  -7.20%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.34%  "

# Row 10 - OKB
Set-TextValue "D10" "58.89"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11 - Avalanche
$ws.Range("E11").Value = "  -10.25%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0936"
$ws.Range("E12").Value = "  -2.95%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.46%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -4.93%  "

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range("D15").Value = "2.505.17"
$ws.Range("E15").Value = "  -1.87%  "

# Row 16 - Chainlink
Set-TextValue "D16" "14.28"
$ws.Range("E16").Value = "  -4.63%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -2.87%  "

# Row 18 - Wrapped Ether
$ws.Range("D18").Value = "2.170.70"
$ws.Range("E18").Value = "  -1.90%  "

# Row 19 - Wrapped BTC
$ws.Range("D19").Value = "41.096.74"
$ws.Range("E19").Value = "  -1.48%  "

# Row 20 - Shiba Inu (price contains subscript-3 character U+2083)
$sub3 = [char]0x2083
$d20Value = "0.0$sub3" + "0945"
$ws.Range("D20").Value = $d20Value
$ws.Range("E20").Value = "  -1.83%  "

# Row 21 - Litecoin
Set-TextValue "D21" "71.63"
$ws.Range("E21").Value = "  -1.70%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.75%  "

# Row 23 - Bitcoin Cash
Set-TextValue "D23" "229.92"

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -2.10%  "

# Row 25 - WEMIX Token
$ws.Range("E25").Value = "  -4.67%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.20%  "

# Row 27 - Cosmos
Set-TextValue "D27" "11.30"
$ws.Range("E27").Value = "  +3.89%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "2.41"
$ws.Range("E28").Value = "  -5.09%  "

# Row 29 - Monero
Set-TextValue "D29" "168.04"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -3.72%  "

# Row 31 - Ethereum Classic
Set-TextValue "D31" "20.17"
$ws.Range("E31").Value = "  -2.98%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.79%  "

# Row 33 - Internet Computer (DFINITY)
Set-TextValue "D33" "5.65"
$ws.Range("E33").Value = "  +1.18%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0746"
$ws.Range("E34").Value = "  +0.87%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  -2.88%  "

# Row 36 - Filecoin
Set-TextValue "D36" "4.51"
$ws.Range("E36").Value = "  -4.50%  "

# Row 37 - Render Token
$ws.Range("E37").Value = "  -1.28%  "

# Row 38 - Injective Protocol
Set-TextValue "D38" "24.64"
$ws.Range("E38").Value = "  -6.59%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0304"
$ws.Range("E39").Value = "  +0.50%  "

# Row 40 - FTX Token
Set-TextValue "D40" "5.49"
$ws.Range("E40").Value = "  +12.48%  "

# Row 41 - Lido DAO Token
$ws.Range("E41").Value = "  -3.78%  "

# Row 42 - THORChain
$ws.Range("E42").Value = "  -6.95%  "

# Row 43 - MultiversX
Set-TextValue "D43" "60.74"
$ws.Range("E43").Value = "  -7.52%  "

# Row 44 - Celestia
Set-TextValue "D44" "11.27"
$ws.Range("E44").Value = "  -8.42%  "

# Row 45 - Frax Share
Set-TextValue "D45" "8.48"
$ws.Range("E45").Value = "  -2.66%  "

# Row 46 - Cronos
Set-TextValue "D46" "0.0991"
$ws.Range("E46").Value = "  -3.21%  "

# Row 47 - Algorand
Set-TextValue "D47" "0.188"
$ws.Range("E47").Value = "  -7.75%  "

# Row 48 - BinanceUSD
$ws.Range("E48").Value = "  -0.11%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  -2.70%  "

# Row 50 - Synthetix Network
Set-TextValue "D50" "4.26"
$ws.Range("E50").Value = "  -10.18%  "

# Row 51 - Trust Wallet Token
$ws.Range("E51").Value = "  -3.86%  "
